$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 86 (shifts existing rows 86..199 down to 87..200)
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the inserted record
$ws.Range("A86").Value = 5
$ws.Range("B86").Value = "Macroferia Regional de Talca"
$ws.Range("C86").Value = "Maule"
$ws.Range("D86").Value = 44467
$ws.Range("E86").Value = 7
$ws.Range("F86").Value = 100112032
$ws.Range("G86").Value = "Zapallo italiano"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = 13000
$ws.Range("N86").Value = "`$/caja 50 unidades"
$ws.Range("O86").Value = "Región de Arica y Parinacota"
$ws.Range("P86").Value = 260
$ws.Range("Q86").Value = 50
$ws.Range("R86").Value = "Hortaliza"
